$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.395.57"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.860.80"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4647"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3720"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07350"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8849"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07912"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("D13").Value = "1.927.06"
$ws.Range("E13").Value = "  +9.98%  "
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.585"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008876"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "27.434.36"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "2.151.08"
$ws.Range("E24").Value = "  +6.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.899"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.077"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.133"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08895"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7579"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.35%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.023"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.162"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.490"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.659"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01968"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.077"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05249"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.113"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5175"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1648"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.354"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4843"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.54%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.653"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.88%  "
